{"js": "// Applies the textual corrections made to the \"Impedimentos\" paragraph\n// that discusses the class diagram / dynamic DB generation tool:\n//\n//   \"... datos asi lo permite , esta tecnichal task debe ser revisada\n//    en el siguiente Sprint como se aclara en en Sprint Review.\"\n// ->\n//   \"... datos as\u00ed lo permite, esta Tecnichal Task debe ser revisada\n//    en pr\u00f3ximos Sprint como se aclara en el Sprint Review.\"\n//\n// Helper: find the first match for `searchText` in the document body and\n// replace it in place with `replacement`, preserving the run's original\n// formatting (insertText onto a search-result range keeps the run's\n// rPr/bold/lang properties).\nasync function replaceFirst(searchText, replacement, matchWholeWord) {\n  const results = context.document.body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: !!matchWholeWord,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"replaceFirst: no match found for \" + JSON.stringify(searchText));\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"asi\" -> \"as\u00ed\" (missing accent)\nawait replaceFirst(\"asi\", \"as\u00ed\", true);\n\n// 2) \"permite ,\" -> \"permite,\" (drop the stray space before the comma)\nawait replaceFirst(\"permite ,\", \"permite,\");\n\n// 3) \"tecnichal\" -> \"Tecnichal\" (capitalised, part of \"Tecnichal Task\")\nawait replaceFirst(\"tecnichal\", \"Tecnichal\", true);\n\n// 4) \"task\" -> \"Task\" (capitalised, part of \"Tecnichal Task\")\nawait replaceFirst(\"task\", \"Task\", true);\n\n// 5) \"el siguiente Sprint\" -> \"pr\u00f3ximos Sprint\"\nawait replaceFirst(\"el siguiente Sprint\", \"pr\u00f3ximos Sprint\");\n\n// 6) \"aclara en en Sprint\" -> \"aclara en el Sprint\" (duplicated \"en\" typo).\n// The duplicated \"en\" sits in its own run wrapped by a <w:proofErr> pair,\n// so rather than sweeping the whole phrase (which would delete that run\n// and orphan its proofErr markers) we isolate just that word: locate the\n// unique \"en Sprint\" match and split it on the space to get an \"en \"\n// sub-range, then replace only that narrow range.\nconst enSprintMatches = context.document.body.search(\"en Sprint\", { matchCase: true });\nenSprintMatches.load(\"items\");\nawait context.sync();\nif (enSprintMatches.items.length === 0) {\n  throw new Error('replaceFirst: no match found for \"en Sprint\"');\n}\nconst enSprintParts = enSprintMatches.items[0].split([\" \"], false, false);\nenSprintParts.load(\"items\");\nawait context.sync();\nenSprintParts.items[0].insertText(\"el \", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Applies the textual corrections made to the \"Impedimentos\" paragraph\n# that discusses the class diagram / dynamic DB generation tool:\n#\n#   \"... datos asi lo permite , esta tecnichal task debe ser revisada\n#    en el siguiente Sprint como se aclara en en Sprint Review.\"\n# ->\n#   \"... datos as\u00ed lo permite, esta Tecnichal Task debe ser revisada\n#    en pr\u00f3ximos Sprint como se aclara en el Sprint Review.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($searchText, $replaceText, $wholeWord) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0            # wdFindStop - don't wrap past the end\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = [bool]$wholeWord\n    $find.MatchWildcards = $false\n    # Replace:=2 -> wdReplaceOne (replace only the first/next match)\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\n# 1) \"asi\" -> \"as\u00ed\" (missing accent)\nReplace-FirstMatch \"asi\" \"as\u00ed\" $true\n\n# 2) \"permite ,\" -> \"permite,\" (drop the stray space before the comma)\nReplace-FirstMatch \"permite ,\" \"permite,\" $false\n\n# 3) \"tecnichal\" -> \"Tecnichal\" (capitalised, part of \"Tecnichal Task\")\nReplace-FirstMatch \"tecnichal\" \"Tecnichal\" $true\n\n# 4) \"task\" -> \"Task\" (capitalised, part of \"Tecnichal Task\")\nReplace-FirstMatch \"task\" \"Task\" $true\n\n# 5) \"el siguiente Sprint\" -> \"pr\u00f3ximos Sprint\"\nReplace-FirstMatch \"el siguiente Sprint\" \"pr\u00f3ximos Sprint\" $false\n\n# 6) \"aclara en en Sprint\" -> \"aclara en el Sprint\" (duplicated \"en\" typo)\nReplace-FirstMatch \"aclara en en Sprint\" \"aclara en el Sprint\" $false\n"}
